$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values for the columns that move: D, J, K, L, M, P
# Rows involved: 2-10 and 12-20 (row 11 is unchanged)
$rows = @(2,3,4,5,6,7,8,9,10,12,13,14,15,16,17,18,19,20)
$cols = @("D","J","K","L","M","P")

$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Mapping: destination row -> source row (values to copy from source row's old snapshot into destination row)
$mapping = @{
    2  = 6
    3  = 19
    4  = 2
    5  = 9
    6  = 14
    7  = 3
    8  = 12
    9  = 8
    10 = 17
    12 = 18
    13 = 20
    14 = 15
    15 = 4
    16 = 7
    17 = 5
    18 = 10
    19 = 13
    20 = 16
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcData[$c]
    }
}
